$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Title"
$ws.Range("B2").Value = "Python"
$ws.Range("C1").Value = "Description"
$ws.Range("C2").Value = "Python Lesson"

$ws.Range("C2").Select()
